$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Read current (pre-edit) values for columns D,L,M,N,O,P,R,S for rows 2..41
$orig = @{}
for ($r = 2; $r -le 41; $r++) {
    $row = @{}
    $row["D"] = $ws.Range("D$r").Value2
    $row["L"] = $ws.Range("L$r").Value2
    $row["M"] = $ws.Range("M$r").Value2
    $row["N"] = $ws.Range("N$r").Value2
    $row["O"] = $ws.Range("O$r").Value2
    $row["P"] = $ws.Range("P$r").Value2
    $row["R"] = $ws.Range("R$r").Value2
    $row["S"] = $ws.Range("S$r").Value2
    $orig[$r] = $row
}

# Mapping: target row -> source row (pre-edit) from which D,L,M,N,O,P,R,S are copied
$map = @{}
$map[2] = 17
$map[3] = 19
$map[4] = 29
$map[5] = 27
$map[6] = 34
$map[7] = 6
$map[8] = 24
$map[9] = 32
$map[10] = 7
$map[11] = 22
$map[12] = 30
$map[13] = 25
$map[14] = 2
$map[15] = 18
$map[16] = 20
$map[17] = 8
$map[18] = 31
$map[19] = 35
$map[20] = 14
$map[21] = 40
$map[22] = 11
$map[23] = 33
$map[24] = 23
$map[25] = 41
$map[26] = 21
$map[27] = 13
$map[28] = 12
$map[29] = 39
$map[30] = 16
$map[31] = 9
$map[32] = 38
$map[33] = 10
$map[34] = 26
$map[35] = 15
$map[36] = 3
$map[37] = 5
$map[38] = 37
$map[39] = 36
$map[40] = 4
$map[41] = 28

foreach ($target in ($map.Keys | Sort-Object)) {
    $source = $map[$target]
    $src = $orig[$source]
    $ws.Range("D$target").Value2 = $src["D"]
    $ws.Range("L$target").Value2 = $src["L"]
    $ws.Range("M$target").Value2 = $src["M"]
    $ws.Range("N$target").Value2 = $src["N"]
    $ws.Range("O$target").Value2 = $src["O"]
    $ws.Range("P$target").Value2 = $src["P"]
    $ws.Range("R$target").Value2 = $src["R"]
    $ws.Range("S$target").Value2 = $src["S"]
}

Write-Host "Row permutation applied."